$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.445.51"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.377.36"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "315.95"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "108.94"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "40.99"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "8.56"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "0.986"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").Value = "2.737.19"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "15.49"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "2.372.10"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "45.400.98"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "16.15"
$ws.Range("E19").Value = "  +21.63%  "
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "0.0000107"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "3.65"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "73.35"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "261.35"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").Value = "11.20"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "0.0967"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "22.40"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "37.32"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").Value = "166.85"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").Value = "2.87"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").Value = "4.08"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("D40").Value = "0.0356"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").Value = "98.07"
$ws.Range("E42").Value = "  -7.37%  "
$ws.Range("D43").Value = "70.60"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "13.16"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("D45").Value = "0.229"
$ws.Range("E45").Value = "  -5.14%  "
$ws.Range("D46").Value = "6.04"
$ws.Range("E46").Value = "  +4.23%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "1.829.40"
$ws.Range("E48").Value = "  +10.72%  "
$ws.Range("D49").Value = "84.29"
$ws.Range("E49").Value = "  +7.41%  "
$ws.Range("D50").Value = "112.09"
$ws.Range("E50").Value = "  -4.61%  "
$ws.Range("D51").Value = "9.28"
$ws.Range("E51").Value = "  -1.03%  "
